# 自动更新Excel文件
# For each data row, the "剩余" (remaining days) counter in column E is
# decremented by one day. When the counter reaches 1 it wraps back to 10
# (the billing/delivery cycle length) and the "开始时间" date in column F
# is advanced by the 10 days of the cycle that just elapsed. Rows whose F
# value is not a well-formed 8-digit yyyyMMdd date are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $eCell = $ws.Cells.Item($row, 5)   # column E
    $fCell = $ws.Cells.Item($row, 6)   # column F

    $eRaw = $eCell.Value2
    $fRaw = $fCell.Value2

    if ($null -eq $eRaw -or $eRaw -eq "") {
        continue
    }

    $eVal = [int]$eRaw
    $fText = "$fRaw"

    $validDate = $false
    if ($fText.Length -eq 8) {
        $year = [int]$fText.Substring(0, 4)
        $month = [int]$fText.Substring(4, 2)
        $day = [int]$fText.Substring(6, 2)
        if ($month -ge 1 -and $month -le 12 -and $day -ge 1 -and $day -le 31) {
            $validDate = $true
        }
    }

    if (-not $validDate) {
        continue
    }

    if ($eVal -eq 1) {
        $newE = 10
        $fDate = [datetime]::ParseExact($fText, "yyyyMMdd", [System.Globalization.CultureInfo]::InvariantCulture)
        $newFDate = $fDate.AddDays(10)
        $newF = [int]$newFDate.ToString("yyyyMMdd")

        $eCell.Value2 = $newE
        $fCell.Value2 = $newF
    }
    else {
        $newE = $eVal - 1
        $eCell.Value2 = $newE
    }
}
